$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setup")

# Row 5: rename the "DATA" field to "ID" and flip its key flag to 0
$ws.Range("B5").Value = "ID"
$ws.Range("C5").Value = 0

# Row 6: add a new "XREF" field (key flag 1, length 35, justification "L")
$ws.Range("B6").Value = "XREF"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 35
$ws.Range("E6").Value = "L"

# Update the active selection to match the authored state
$ws.Range("E7").Select()
